$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# Overview: Latest HO Xliff Generate Date (G2)
$wsOverview.Range("G2").Value = "2016-10-21 00:56:17"

# zh-cn: Correspond Handoff Datetime (H2), Correspond Handback DateTime (K2)
$wsZhCn.Range("H2").Value = "2016-10-21 00:56:05"
$wsZhCn.Range("K2").Value = "2016-10-21 00:56:48"

# de-de: Correspond Handoff Datetime (H2), Correspond Handback DateTime (K2)
$wsDeDe.Range("H2").Value = "2016-10-21 00:56:17"
$wsDeDe.Range("K2").Value = "2016-10-21 00:57:06"
